$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 37; this shifts the existing rows 37-42 down to 38-43.
$ws.Rows.Item(37).Insert()

# The new row 37 is a new price record for Albahaca, same market/region/etc. as
# its neighbours, copy the constant columns from the row right below it (the
# old row 37, now shifted to row 38).
$ws.Range("A37").Value = 1
$ws.Range("B37").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C37").Value = "Arica y Parinacota"
$ws.Range("D37").Value = 44798
$ws.Range("E37").Value = 15
$ws.Range("F37").Value = 100112052
$ws.Range("G37").Value = "Albahaca"
$ws.Range("H37").Value = "Sin especificar"
$ws.Range("I37").Value = "Primera"
$ws.Range("J37").Value = 250
$ws.Range("K37").Value = 1500
$ws.Range("L37").Value = 2000
$ws.Range("M37").Value = 1750
$ws.Range("N37").Value = "$/paquete"
$ws.Range("O37").Value = "Región de Arica y Parinacota"
$ws.Range("P37").Value = 1750
$ws.Range("Q37").Value = 1
$ws.Range("R37").Value = "Hortaliza"

# Match the date cell's number format to the other date cells in column D.
$ws.Range("D37").NumberFormat = $ws.Range("D38").NumberFormat
